# Update the dSF column (F) with repulled / recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = -3
    4  = -1
    5  = -5
    7  = 2
    8  = -4
    9  = -4
    11 = -1
    12 = 3
    14 = -1
    15 = -4
    16 = 0
    17 = -10
    19 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
